# The document body is a single paragraph consisting mostly of anchored
# (floating) drawings plus one trailing stray run containing the single
# character "c". The commit removes that leftover run entirely.
#
# Use Find/Replace on the document's main story so the run (together with
# its run properties) is deleted rather than merely emptied.
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "c",     # FindText
    $false,  # MatchCase
    $false,  # MatchWholeWord
    $false,  # MatchWildcards
    $false,  # MatchSoundsLike
    $false,  # MatchAllWordForms
    $true,   # Forward
    1,       # Wrap (wdFindContinue)
    $false,  # Format
    "",      # ReplaceWith
    2        # Replace (wdReplaceAll)
)
